$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SC 92" row (row 28) entirely; remaining rows below shift up by one.
$ws.Rows.Item(28).Delete()

# Remove the "RM 232" row (row 26) entirely; remaining rows below shift up by one.
$ws.Rows.Item(26).Delete()

# After the two deletions, rows 26-33 now hold (in order):
#   26 SC 5, 27 SC 101, 28 SC 105, 29 SC 119, 30 SC 120, 31 SC 132, 32 SC 193, 33 SC 232
# Apply the remaining value corrections to match the target data.

# SC 5 row: fill in column B (was blank)
$ws.Cells.Item(26, 2).Value = -20.2

# SC 101 row: clear column B (becomes blank)
$ws.Cells.Item(27, 2).Value = ""

# SC 232 row: fill in column C (was blank)
$ws.Cells.Item(33, 3).Value = 10.4
